$wb = $excel.ActiveWorkbook
Write-Output $wb.Worksheets.Count
foreach ($ws in $wb.Worksheets) {
    Write-Output $ws.Name
}
